# Update Mappings 22 Ontologies
# Adds a new "EDAM_DEF" column (F) to the CAO_EDAM mapping sheet, containing
# the EDAM definition text for each mapped row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Cells.Item(1, 6).Value = "EDAM_DEF"

# Copy the header formatting (bold font, border, centered/top alignment)
# from the neighboring EDAM_DESC header cell onto the new EDAM_DEF header cell.
$ws.Cells.Item(1, 5).Copy() | Out-Null
$ws.Cells.Item(1, 6).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# EDAM_DEF values for each data row (rows 2-13), matching the CAO/EDAM mapping
$edamDefs = @(
    "['An array of numerical values.']",
    "['A specification of a chemical structure in SMILES format.']",
    "['The InChIKey (hashed InChI) is a fixed length (25 character) condensed digital representation of an InChI chemical structure specification. It uniquely identifies a chemical compound.']",
    "['Chemical structure specified in IUPAC International Chemical Identifier (InChI) line notation.']",
    "['An array of numerical values.']",
    "['Unique identifier of a chemical compound.']",
    "['A human-readable collection of information about about how a scientific experiment or analysis was carried out that results in a specific set of data or results used for further analysis or to test a specific hypothesis.']",
    "['The concentration of a chemical compound.']",
    "['The spectrum of frequencies of electromagnetic radiation emitted from a molecule as a result of some spectroscopy experiment.']",
    "['A report about localisation of the isolaton of biological material e.g. country or coordinates.']",
    "['An alignment of molecular sequences, structures or profiles derived from them.']",
    "['A valid email address of an end-user.']"
)

for ($i = 0; $i -lt $edamDefs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $edamDefs[$i]
}
